# Updated cryptos list on Tue Apr 23 07:35:17 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to write text-typed values into the D (Price) column
# without Excel auto-converting numeric-looking strings (e.g. "601.26")
# into real numbers, and without leaving any NumberFormat/style residue
# on the destination cells (copy / paste-special values-only keeps the
# destination formatting untouched).
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

function Set-TextValue($targetAddr, $val) {
    $scratch.Value = $val
    $scratch.Copy()
    $ws.Range($targetAddr).PasteSpecial(-4163)
}

Set-TextValue "D2" "66.272.21"
$ws.Range("E2").Value = "  +0.33%  "
Set-TextValue "D3" "3.163.22"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue "D5" "601.26"
$ws.Range("E5").Value = "  -0.50%  "
Set-TextValue "D6" "153.58"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +2.75%  "
Set-TextValue "D9" "3.158.73"
$ws.Range("E9").Value = "  -1.56%  "
Set-TextValue "D10" "0.157"
$ws.Range("E10").Value = "  -2.48%  "
$ws.Range("E11").Value = "  -9.75%  "
Set-TextValue "D12" "0.508"
$ws.Range("E12").Value = "  -0.12%  "
Set-TextValue "D13" "0.0000266"
$ws.Range("E13").Value = "  -3.16%  "
Set-TextValue "D14" "38.30"
$ws.Range("E14").Value = "  -1.16%  "
Set-TextValue "D15" "3.686.47"
$ws.Range("E15").Value = "  -1.29%  "
Set-TextValue "D16" "66.322.24"
$ws.Range("E16").Value = "  +0.18%  "
Set-TextValue "D17" "7.35"
$ws.Range("E17").Value = "  -1.39%  "
Set-TextValue "D18" "3.166.89"
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("E19").Value = "  +0.19%  "
Set-TextValue "D20" "509.22"
$ws.Range("E20").Value = "  -0.44%  "
Set-TextValue "D21" "15.32"
$ws.Range("E21").Value = "  -1.86%  "
Set-TextValue "D22" "0.726"
$ws.Range("E22").Value = "  -1.41%  "
Set-TextValue "D23" "8.06"
$ws.Range("E23").Value = "  +0.72%  "
Set-TextValue "D24" "14.63"
$ws.Range("E24").Value = "  -3.99%  "
Set-TextValue "D25" "84.50"
$ws.Range("E25").Value = "  -0.97%  "
Set-TextValue "D26" "0.998"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -1.29%  "
Set-TextValue "D28" "9.08"
$ws.Range("E28").Value = "  -1.91%  "
Set-TextValue "D29" "2.37"
$ws.Range("E29").Value = "  +5.80%  "
Set-TextValue "D30" "3.04"
$ws.Range("E30").Value = "  +5.59%  "
Set-TextValue "D31" "6.86"
$ws.Range("E31").Value = "  -0.48%  "
Set-TextValue "D32" "27.82"
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  -1.58%  "
Set-TextValue "D35" "6.47"
$ws.Range("E35").Value = "  -2.59%  "
Set-TextValue "D36" "509.48"
$ws.Range("E36").Value = "  +5.07%  "
Set-TextValue "D37" "54.74"
$ws.Range("E37").Value = "  -1.13%  "
Set-TextValue "D38" "0.0884"
$ws.Range("E38").Value = "  -3.26%  "
Set-TextValue "D39" "0.0417"
$ws.Range("E39").Value = "  -0.94%  "
Set-TextValue "D40" "0.127"
$ws.Range("E40").Value = "  +5.87%  "
Set-TextValue "D41" "8.78"
$ws.Range("E41").Value = "  -0.74%  "
Set-TextValue "D42" "0.0₃0673"
$ws.Range("E42").Value = "  +4.24%  "
Set-TextValue "D43" "0.295"
$ws.Range("E43").Value = "  -0.53%  "
Set-TextValue "D44" "2.77"
$ws.Range("E44").Value = "  -7.49%  "
Set-TextValue "D45" "2.41"
$ws.Range("E45").Value = "  -3.66%  "
Set-TextValue "D46" "2.825.75"
$ws.Range("E46").Value = "  -4.41%  "
Set-TextValue "D47" "27.82"
$ws.Range("E47").Value = "  -3.76%  "
Set-TextValue "D49" "2.35"
$ws.Range("E49").Value = "  +0.81%  "
Set-TextValue "D50" "0.116"
$ws.Range("E50").Value = "  +0.18%  "
Set-TextValue "D51" "34.64"
$ws.Range("E51").Value = "  +1.39%  "

# Clean up the scratch cell so it leaves no trace in the saved sheet.
$scratch.Clear()
$excel.CutCopyMode = $false
